$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '97.758.18'
$ws.Range('E2').Value = '  +0.31%  '

# Row 3
$ws.Range('D3').Value = '3.302.14'
$ws.Range('E3').Value = '  -1.54%  '

# Row 4
$ws.Range('E4').Value = '  +0.07%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '255.58'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.36%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '620.82'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.77%  '

# Row 7
$ws.Range('E7').Value = '  +25.87%  '

# Row 8
$ws.Range('E8').Value = '  +2.36%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('D9').Style = 'Normal'

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.898'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +13.38%  '

# Row 11
$ws.Range('D11').Value = '3.299.20'
$ws.Range('E11').Value = '  -1.79%  '

# Row 12
$ws.Range('E12').Value = '  -0.79%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '38.62'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +9.52%  '

# Row 14
$ws.Range('D14').Value = '97.464.47'
$ws.Range('E14').Value = '  +0.49%  '

# Row 15
$ws.Range('E15').Value = '  -0.62%  '

# Row 16
$ws.Range('D16').Value = '3.919.96'
$ws.Range('E16').Value = '  -1.10%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.45'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.39%  '

# Row 18
$ws.Range('D18').Value = '3.303.40'
$ws.Range('E18').Value = '  -1.06%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.52'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.75%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.10'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.77%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.14'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.55%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '477.61'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.76%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.41'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.78%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000204'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.15%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.58'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.33%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '87.69'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.85%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.83'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.80%  '

# Row 28
$ws.Range('D28').Value = '3.473.48'
$ws.Range('E28').Value = '  -1.81%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.294'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +20.96%  '

# Row 30
$ws.Range('E30').Value = '  -0.12%  '

# Row 31
$ws.Range('E31').Value = '  +2.79%  '

# Row 32
$ws.Range('E32').Value = '  +7.41%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '9.82'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.93%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.12%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '27.49'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.85%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '7.15'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.93%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.147'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.08%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.93'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.34%  '

# Row 39
$ws.Range('E39').Value = '  +0.61%  '

# Row 40
$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '489.52'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.01%  '

# Row 41
$ws.Range('B41').Value = 'PolygonEcosystemToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.453'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.30%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.67'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.37%  '

# Row 43
$ws.Range('E43').Value = '  -3.47%  '

# Row 44
$ws.Range('E44').Value = '  -1.88%  '

# Row 45
$ws.Range('E45').Value = '  -0.01%  '

# Row 46
$ws.Range('E46').Value = '  -4.54%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '157.86'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.18%  '

# Row 48
$ws.Range('E48').Value = '  -2.46%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.839'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.12%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.63'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.21%  '

# Row 51
$ws.Range('B51').Value = 'Cosmos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.07'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +10.66%  '
